$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect List")
$ws.Activate()

# --- New defect entry: row 60 -------------------------------------------
# Copy existing cell formatting (format-only paste) onto the new row so the
# workbook's shared style table is reused rather than growing.
$ws.Range("A8").Copy()
$ws.Range("A60").PasteSpecial(-4122)

$ws.Range("B52").Copy()
$ws.Range("B60").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C60").PasteSpecial(-4122)

$ws.Range("D52").Copy()
$ws.Range("D60").PasteSpecial(-4122)

$ws.Range("E52").Copy()
$ws.Range("E60").PasteSpecial(-4122)

$ws.Range("F52").Copy()
$ws.Range("F60").PasteSpecial(-4122)

$ws.Range("G52").Copy()
$ws.Range("G60").PasteSpecial(-4122)

$ws.Range("H8").Copy()
$ws.Range("H60").PasteSpecial(-4122)

$ws.Range("I52").Copy()
$ws.Range("I60").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Values -- set G/H/F in this order so any newly introduced shared strings
# land in the table in the same order the authoring session produced them.
$ws.Range("A60").Value = 53
$ws.Range("B60").Value = "Steven "
$ws.Range("C60").Value = Get-Date -Year 2016 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0
$ws.Range("D60").Value = 4
$ws.Range("E60").Value = "Screen "
$ws.Range("G60").Value = "Main Menu button now takes user to Theme"
$ws.Range("H60").Value = "Click the main menu button"
$ws.Range("F60").Value = "Main Menu Button"
$ws.Range("I60").Value = "Fixed"

$ws.Rows.Item(60).RowHeight = 31.5

# --- View state: scrolled down to show the new row, selection moved ------
$ws.Range("F62").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
